$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old "Interval time" data (rows 1-5) before writing the new
# "Phase" data (rows 1-3 only).
$ws.Range("A1:B5").ClearContents()

# Column A: phase labels (written first so their shared-string entries
# come before the column B values).
$ws.Cells.Item(1, 1).Value = "Phase 1: "
$ws.Cells.Item(2, 1).Value = "Phase 2: "
$ws.Cells.Item(3, 1).Value = "Phase 3: "

# Column B: TF lists for each phase.
$ws.Cells.Item(1, 2).Value = "TF13, TF16, TF1, TF2, TF3, TF4, TF5, TF6, TF7"
$ws.Cells.Item(2, 2).Value = "TF13, TF16, TF3, TF4, TF5, TF6, TF7, TF8, TF9, TF11, TF14, TF15"
$ws.Cells.Item(3, 2).Value = "TF13, TF16, TF3, TF6, TF10, TF11, TF12, TF14, TF15"

# Widen column B to fit the longer TF list text.
$ws.Columns.Item(2).ColumnWidth = 24.666666666666668

# Move the active selection to just past the new data (one row below).
$ws.Range("B4").Select()
